$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
Set-TextValue 2 4 "37.368.88"
Set-TextValue 2 5 "  -1.26%  "

# Row 3
Set-TextValue 3 4 "2.052.59"
Set-TextValue 3 5 "  -1.45%  "

# Row 4
Set-TextValue 4 5 "  -0.06%  "

# Row 5
Set-TextValue 5 4 "230.74"
Set-TextValue 5 5 "  -1.02%  "

# Row 7
Set-TextValue 7 5 "  +0.02%  "

# Row 8
Set-TextValue 8 5 "  -3.76%  "

# Row 9
Set-TextValue 9 5 "  -2.79%  "

# Row 10
Set-TextValue 10 4 "0.0772"
Set-TextValue 10 5 "  -2.19%  "

# Row 11
Set-TextValue 11 5 "  +1.40%  "

# Row 12
Set-TextValue 12 4 "2.353.23"
Set-TextValue 12 5 "  -1.53%  "

# Row 13
Set-TextValue 13 4 "14.61"
Set-TextValue 13 5 "  -0.98%  "

# Row 14
Set-TextValue 14 4 "20.67"
Set-TextValue 14 5 "  -2.48%  "

# Row 15
Set-TextValue 15 4 "0.757"
Set-TextValue 15 5 "  -2.49%  "

# Row 16
Set-TextValue 16 4 "5.28"
Set-TextValue 16 5 "  -1.29%  "

# Row 17
Set-TextValue 17 4 "2.058.52"
Set-TextValue 17 5 "  -0.73%  "

# Row 18
Set-TextValue 18 4 "37.287.99"
Set-TextValue 18 5 "  -1.26%  "

# Row 19
Set-TextValue 19 4 "6.09"
Set-TextValue 19 5 "  -0.96%  "

# Row 20
Set-TextValue 20 4 "69.61"
Set-TextValue 20 5 "  -2.74%  "

# Row 21
Set-TextValue 21 4 "0.0₃0823"
Set-TextValue 21 5 "  -3.21%  "

# Row 22
Set-TextValue 22 4 "226.35"
Set-TextValue 22 5 "  -0.78%  "

# Row 24
Set-TextValue 24 4 "2.39"
Set-TextValue 24 5 "  +0.05%  "

# Row 25
Set-TextValue 25 5 "  -3.57%  "

# Row 26
Set-TextValue 26 4 "9.89"
Set-TextValue 26 5 "  +8.11%  "

# Row 27
Set-TextValue 27 4 "169.98"
Set-TextValue 27 5 "  -0.94%  "

# Row 28
Set-TextValue 28 5 "  -6.61%  "

# Row 29
Set-TextValue 29 5 "  -1.42%  "

# Row 30
Set-TextValue 30 5 "  -5.11%  "

# Row 31
Set-TextValue 31 5 "  +0.16%  "

# Row 32
Set-TextValue 32 5 "  -4.18%  "

# Row 33
Set-TextValue 33 5 "  -1.50%  "

# Row 34
Set-TextValue 34 5 "  -3.38%  "

# Row 35
Set-TextValue 35 5 "  -1.45%  "

# Row 36
Set-TextValue 36 5 "  -0.04%  "

# Row 37
Set-TextValue 37 5 "  -4.48%  "

# Row 38
Set-TextValue 38 4 "1.00"
Set-TextValue 38 5 "  +0.16%  "

# Row 40
Set-TextValue 40 5 "  +3.11%  "

# Row 41
Set-TextValue 41 4 "98.25"
Set-TextValue 41 5 "  -0.96%  "

# Row 42
Set-TextValue 42 5 "  -2.95%  "

# Row 43
Set-TextValue 43 4 "2.89"
Set-TextValue 43 5 "  -0.09%  "

# Row 44
Set-TextValue 44 4 "1.477.32"
Set-TextValue 44 5 "  +2.16%  "

# Row 45
Set-TextValue 45 5 "  +2.68%  "

# Row 46
Set-TextValue 46 4 "16.62"
Set-TextValue 46 5 "  -0.72%  "

# Row 47
Set-TextValue 47 5 "  -3.04%  "

# Row 48
Set-TextValue 48 2 "FraxShare"
Set-TextValue 48 3 "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue 48 4 "7.27"
Set-TextValue 48 5 "  -1.67%  "

# Row 49
Set-TextValue 49 2 "FTXToken"
Set-TextValue 49 3 "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue 49 4 "3.95"
Set-TextValue 49 5 "  -5.19%  "

# Row 50
Set-TextValue 50 5 "  -1.78%  "

# Row 51
Set-TextValue 51 4 "2.240.16"
Set-TextValue 51 5 "  -1.50%  "

Write-Output "Applied all crypto list updates"